# Fruta / hortaliza, semanal
#
# Inserts a new week of Chirimoya (Cultivar IV Region, Provincia del Elqui)
# price observations dated 2021-09-30 (Excel serial 44469) as rows 18-23,
# pushing the existing rows 18-53 down to rows 24-59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18 downward by 6 to make room for the new weekly block.
$ws.Rows("18:23").Insert()

# Columns that stay constant for every Chirimoya / Cultivar IV Region row.
$A = 9
$B = "Vega Central Mapocho de Santiago"
$C = "Metropolitana"
$E = 13
$F = "Fruta"
$G = 100107
$H = "Otros"
$I = 100107002
$J = "Chirimoya"
$K = "Cultivar IV Región"

$newDate = 44469

function Set-ChirimoyaRow($Row, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $newDate
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

Set-ChirimoyaRow 18 "Cuarta"                  200 1000  1000  1000  "`$/kilo (en caja de 15 kilos)" "Provincia del Elquí" 1000 1
Set-ChirimoyaRow 19 "Especial"                220 20000 20000 20000 "`$/bandeja 8 kilos"            "Provincia del Elquí" 2500 8
Set-ChirimoyaRow 20 "Extra (doble especial)"  180 24000 24000 24000 "`$/bandeja 8 kilos"            "Provincia del Elquí" 3000 8
Set-ChirimoyaRow 21 "Primera"                 280 16000 16000 16000 "`$/bandeja 8 kilos"            "Provincia del Elquí" 2000 8
Set-ChirimoyaRow 22 "Segunda"                 150 13600 13600 13600 "`$/bandeja 8 kilos"            "Provincia del Elquí" 1700 8
Set-ChirimoyaRow 23 "Tercera"                 180 1400  1400  1400  "`$/kilo (en caja de 15 kilos)" "Provincia del Elquí" 1400 1
